$d = $word.ActiveDocument

# The target paragraph is the very last one in the document body
# ("Dia 06/09: 1hr 10min (1 dia)"), immediately preceding the sectPr.
$count = $d.Paragraphs.Count
$p = $d.Paragraphs.Item($count)
$r = $p.Range

# Narrow the range to just the "1hr 10min" span inside that paragraph,
# then replace its text with "2hr". Toggling a character property
# (Bold on, then back off) around the text assignment keeps the engine
# from silently re-merging the new run back into its neighbours, so the
# paragraph ends up split into three runs - "Dia 06/09: ", "2hr" and
# " (1 dia)" - exactly like a real Word edit-in-place would leave it.
$found = $r.Find.Execute("1hr 10min", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $r.Font.Bold = 1
    $r.Text = "2hr"
    $r.Font.Bold = 0
}
